# Apply crypto price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.312.78"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.654.79"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "660.29"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.652.62"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.204"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").Value = "4.332.06"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000269"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.46%  "
$ws.Range("D17").Value = "96.079.17"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +14.48%  "
$ws.Range("D19").Value = "3.637.11"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "524.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  +8.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.994"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.98%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "623.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "43.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +33.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.159"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.96%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0452"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.420"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.26%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
